# Adding Code For Payer Operation Part_1
# Appends the next PEP patient id to the patient list worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pepID")
if (-not $ws) {
    $ws = $wb.ActiveSheet
}

$ws.Range("A2").Value = "PEP_ID-2009515"
